$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9215428233146667
$ws.Range("B1").Value = 1.400839924812317
$ws.Range("C1").Value = 2.622637271881104
$ws.Range("D1").Value = 1.512464642524719
$ws.Range("E1").Value = 1.430940508842468
